$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row2
$ws.Range('D2').Value = '42.973.56'
$ws.Range('E2').Value = '  +0.06%  '

# row3
$ws.Range('D3').Value = '2.211.77'
$ws.Range('E3').Value = '  -0.91%  '

# row4
$ws.Range('E4').Value = '  -0.06%  '

# row5
$ws.Range('D5').Formula = "'257.62"
$ws.Range('E5').Value = '  +2.65%  '

# row6
$ws.Range('D6').Formula = "'0.619"
$ws.Range('E6').Value = '  +0.91%  '

# row7
$ws.Range('D7').Formula = "'77.18"
$ws.Range('E7').Value = '  +3.04%  '

# row8
$ws.Range('E8').Value = '  -0.06%  '

# row9
$ws.Range('D9').Formula = "'0.595"
$ws.Range('E9').Value = '  -0.01%  '

# row10
$ws.Range('D10').Formula = "'42.82"
$ws.Range('E10').Value = '  +3.79%  '

# row11
$ws.Range('D11').Formula = "'0.0914"
$ws.Range('E11').Value = '  -0.82%  '

# row12
$ws.Range('D12').Formula = "'7.00"
$ws.Range('E12').Value = '  +1.93%  '

# row13
$ws.Range('E13').Value = '  +1.13%  '

# row14
$ws.Range('D14').Value = '2.542.78'
$ws.Range('E14').Value = '  -1.04%  '

# row15
$ws.Range('D15').Formula = "'14.48"
$ws.Range('E15').Value = '  +0.01%  '

# row16
$ws.Range('D16').Value = '2.212.67'
$ws.Range('E16').Value = '  -0.83%  '

# row17
$ws.Range('D17').Formula = "'0.786"
$ws.Range('E17').Value = '  -0.07%  '

# row18
$ws.Range('D18').Value = '42.914.50'
$ws.Range('E18').Value = '  +0.14%  '

# row19
$ws.Range('E19').Value = '  +0.21%  '

# row20
$ws.Range('D20').Formula = "'71.18"
$ws.Range('E20').Value = '  -0.01%  '

# row21
$ws.Range('D21').Formula = "'5.98"
$ws.Range('E21').Value = '  +0.99%  '

# row22
$ws.Range('D22').Formula = "'2.35"
$ws.Range('E22').Value = '  +8.32%  '

# row23
$ws.Range('D23').Formula = "'229.94"
$ws.Range('E23').Value = '  +0.26%  '

# row24
$ws.Range('D24').Formula = "'9.22"
$ws.Range('E24').Value = '  -2.04%  '

# row25
$ws.Range('E25').Value = '  -0.02%  '

# row26
$ws.Range('D26').Formula = "'42.37"
$ws.Range('E26').Value = '  +8.76%  '

# row27
$ws.Range('D27').Formula = "'10.75"
$ws.Range('E27').Value = '  +0.86%  '

# row28
$ws.Range('E28').Value = '  -2.64%  '

# row29
$ws.Range('E29').Value = '  -0.54%  '

# row30
$ws.Range('D30').Formula = "'2.19"
$ws.Range('E30').Value = '  +2.45%  '

# row31
$ws.Range('D31').Formula = "'172.87"
$ws.Range('E31').Value = '  +0.84%  '

# row32
$ws.Range('D32').Formula = "'20.37"
$ws.Range('E32').Value = '  +1.12%  '

# row33
$ws.Range('D33').Formula = "'0.0865"
$ws.Range('E33').Value = '  +9.02%  '

# row34
$ws.Range('D34').Formula = "'5.23"
$ws.Range('E34').Value = '  +0.60%  '

# row35
$ws.Range('D35').Formula = "'0.121"
$ws.Range('E35').Value = '  +0.18%  '

# row36
$ws.Range('E36').Value = '  +11.37%  '

# row37
$ws.Range('E37').Value = '  -4.39%  '

# row38
$ws.Range('E38').Value = '  -0.74%  '

# row39
$ws.Range('D39').Formula = "'12.84"
$ws.Range('E39').Value = '  +3.26%  '

# row40
$ws.Range('E40').Value = '  +20.29%  '

# row41
$ws.Range('D41').Formula = "'2.11"
$ws.Range('E41').Value = '  +1.14%  '

# row42
$ws.Range('B42').Value = 'MultiversX'
$ws.Range('C42').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D42').Formula = "'61.44"
$ws.Range('E42').Value = '  +3.55%  '

# row43
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').Formula = "'0.202"
$ws.Range('E43').Value = '  -0.80%  '

# row44
$ws.Range('D44').Formula = "'5.27"
$ws.Range('E44').Value = '  -1.68%  '

# row45
$ws.Range('D45').Formula = "'103.09"
$ws.Range('E45').Value = '  +0.01%  '

# row46
$ws.Range('D46').Formula = "'8.49"
$ws.Range('E46').Value = '  -1.91%  '

# row47
$ws.Range('D47').Formula = "'0.468"
$ws.Range('E47').Value = '  -2.72%  '

# row48
$ws.Range('E48').Value = '  -1.46%  '

# row49
$ws.Range('E49').Value = '  +0.26%  '

# row50
$ws.Range('E50').Value = '  +0.04%  '

# row51
$ws.Range('D51').Formula = "'1.45"
$ws.Range('E51').Value = '  +21.06%  '

